# edit.ps1 - apply "output generated at 456a3b4" update to 上海-漫展信息.xlsx
# Source workbook is already open as $excel.ActiveWorkbook
$wb = $excel.ActiveWorkbook

# --- Sheet handles ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- sheet1 (展览): refresh row 5's view-count, insert new row 6 event, shift rows 6-9 content into 7-10 ---
$wsExpo.Range("F5").Value = 2837

# Row 6
$wsExpo.Range("C6").Value = "上海·GUTI&GUTA LAND"
$wsExpo.Range("D6").Value = "翔殷路1099号 合生汇"
$wsExpo.Range("E6").Value = "2024.07.26 10:00-08.31 22:00"
$wsExpo.Range("F6").Value = 224
$wsExpo.Range("G6").Value = 34.9
$wsExpo.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=89666"
$wsExpo.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202407/eeb6eAaB1721378157456.png"

# Row 7
$wsExpo.Range("B7").NumberFormat = "@"
$wsExpo.Range("B7").Value = "2024-07-26"
$wsExpo.Range("C7").Value = "上海·英雄的苍穹：正子公也三国、水浒绘画艺术大展"
$wsExpo.Range("D7").Value = "东安路888号 W+艺术中心"
$wsExpo.Range("E7").Value = "2024.07.26 10:00-08.25 19:00"
$wsExpo.Range("F7").Value = 14
$wsExpo.Range("G7").Value = 38
$wsExpo.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=89497"
$wsExpo.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202407/LGSorblv1721204272396.jpeg"

# Row 8
$wsExpo.Range("C8").Value = "上海·THE哆啦A梦展"
$wsExpo.Range("D8").Value = "昭化路638号CPARK.D栋 iag艺术院线"
$wsExpo.Range("E8").Value = "2024.08.02 10:00-10.05 22:00"
$wsExpo.Range("F8").Value = 281
$wsExpo.Range("G8").Value = 80
$wsExpo.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=88298"
$wsExpo.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202406/URdR4tbi1719565900366.jpeg"

# Row 9
$wsExpo.Range("B9").NumberFormat = "@"
$wsExpo.Range("B9").Value = "2024-08-02"
$wsExpo.Range("C9").Value = "上海·艺术与潮流·遇见EVA 中国首展"
$wsExpo.Range("D9").Value = "西藏北路166号 静安大悦城北座"
$wsExpo.Range("E9").Value = "2024.08.02 10:00-10.07 22:00"
$wsExpo.Range("F9").Value = 6583
$wsExpo.Range("G9").Value = 89
$wsExpo.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=89161"
$wsExpo.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202407/z8YTdxA71720679877329.jpeg"

# Row 10
$wsExpo.Range("B10").NumberFormat = "@"
$wsExpo.Range("B10").Value = "2024-08-09"
$wsExpo.Range("C10").Value = "上海·新井里美粉丝见面会"
$wsExpo.Range("D10").Value = "西藏南路1号 上海大世界"
$wsExpo.Range("E10").Value = "2024.08.09 16:30-08.11 16:30"
$wsExpo.Range("F10").Value = 7
$wsExpo.Range("G10").Value = 198
$wsExpo.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=90332"
$wsExpo.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202408/geUVjFXF1722842759283.jpeg"

# --- Remaining independent 想去人数 (F-column) refreshes across all sheets ---
# 展览
$sheet1_updates = @(
    @{Row=11; Value=5104},
    @{Row=12; Value=12},
    @{Row=13; Value=570},
    @{Row=14; Value=2710},
    @{Row=18; Value=333},
    @{Row=19; Value=134},
    @{Row=20; Value=148},
    @{Row=21; Value=1118},
    @{Row=22; Value=267},
    @{Row=23; Value=1408},
    @{Row=25; Value=2130},
    @{Row=26; Value=1354},
    @{Row=28; Value=59},
    @{Row=29; Value=1019},
    @{Row=30; Value=52},
    @{Row=31; Value=129},
    @{Row=32; Value=1552},
    @{Row=33; Value=13},
    @{Row=35; Value=1626},
    @{Row=36; Value=1107},
    @{Row=39; Value=328},
    @{Row=40; Value=2346},
    @{Row=41; Value=2606},
    @{Row=42; Value=60},
    @{Row=43; Value=160},
    @{Row=46; Value=287},
    @{Row=48; Value=124},
    @{Row=49; Value=395}
)
foreach ($u in $sheet1_updates) {
    $wsExpo.Cells.Item($u.Row, 6).Value = $u.Value
}

# 演出
$sheet2_updates = @(
    @{Row=8; Value=337},
    @{Row=10; Value=173},
    @{Row=11; Value=105},
    @{Row=12; Value=205},
    @{Row=16; Value=169},
    @{Row=17; Value=48},
    @{Row=23; Value=6},
    @{Row=26; Value=434},
    @{Row=27; Value=33},
    @{Row=29; Value=18},
    @{Row=40; Value=32}
)
foreach ($u in $sheet2_updates) {
    $wsShow.Cells.Item($u.Row, 6).Value = $u.Value
}

# 本地生活
$sheet3_updates = @(
    @{Row=8; Value=1570},
    @{Row=9; Value=1833},
    @{Row=10; Value=2608},
    @{Row=11; Value=911},
    @{Row=12; Value=802},
    @{Row=14; Value=170}
)
foreach ($u in $sheet3_updates) {
    $wsLocal.Cells.Item($u.Row, 6).Value = $u.Value
}

# 全部类型
$sheet4_updates = @(
    @{Row=6; Value=2837},
    @{Row=7; Value=1570},
    @{Row=8; Value=281},
    @{Row=9; Value=2608},
    @{Row=10; Value=6583},
    @{Row=11; Value=911},
    @{Row=12; Value=802},
    @{Row=13; Value=5104},
    @{Row=14; Value=2710},
    @{Row=18; Value=333},
    @{Row=19; Value=148},
    @{Row=20; Value=337},
    @{Row=21; Value=1118},
    @{Row=22; Value=267},
    @{Row=23; Value=105},
    @{Row=24; Value=170},
    @{Row=25; Value=1408},
    @{Row=27; Value=2130},
    @{Row=28; Value=1354},
    @{Row=30; Value=59},
    @{Row=31; Value=169},
    @{Row=32; Value=1019},
    @{Row=33; Value=52},
    @{Row=34; Value=48},
    @{Row=35; Value=1552},
    @{Row=37; Value=1107},
    @{Row=39; Value=434},
    @{Row=40; Value=328},
    @{Row=41; Value=33},
    @{Row=43; Value=2346},
    @{Row=44; Value=2606},
    @{Row=45; Value=160},
    @{Row=46; Value=287},
    @{Row=48; Value=395}
)
foreach ($u in $sheet4_updates) {
    $wsAll.Cells.Item($u.Row, 6).Value = $u.Value
}

